# "shading properties of Singapore archetypes - added shading type to
# construction archetypes"
#
# The ARCHITECTURE sheet's type_shade column (L), rows 2:19, was previously
# blank; populate it with 0 (no shading) for every construction archetype
# row, and make ARCHITECTURE the active sheet/selection (it previously
# wasn't - INTERNAL_LOADS was the last-active tab).

$wb = $excel.ActiveWorkbook

$wsArch = $wb.Worksheets.Item("ARCHITECTURE")

# type_shade (column L), rows 2-19: blank -> 0
$wsArch.Range("L2:L19").Value = 0

# Make ARCHITECTURE the active sheet and select the newly-filled range,
# matching the updated view/selection state in the workbook.
$wsArch.Activate()
$wsArch.Range("L2:L19").Select()
